# Update the "Pais" COVID-19 dashboard sheet with the latest figures and
# refresh the "last updated" timestamp in A1 (cell A1).
#
# The underlying data is kept sorted by column B (Casos totales) descending,
# so when a country's case total overtakes a neighbouring row's total, the
# two rows swap countries: rows 42/43 (Filipinas overtakes Suiza) and rows
# 113/114 (Madagascar overtakes Guinea Ecuatorial) each get a new country
# name in column A in addition to refreshed figures; every other changed
# row keeps its existing country and just gets refreshed figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Last updated" footer, cell A1
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 23 de Junio de 2020 a las 13:25"

# Columns: 1=Pais 2=Casos totales 3=Nuevos casos 4=Casos activos
#          5=Recuperados 6=Casos criticos 7=Muertes hoy 8=Muertes

# Row 7: India
$ws.Cells.Item(7, 2).Value = 441924
$ws.Cells.Item(7, 3).Value = 1474
$ws.Cells.Item(7, 4).Value = 248629
$ws.Cells.Item(7, 5).Value = 179267
$ws.Cells.Item(7, 7).Value = 13
$ws.Cells.Item(7, 8).Value = 14028

# Row 13: Iran
$ws.Cells.Item(13, 2).Value = 209970
$ws.Cells.Item(13, 3).Value = 2445
$ws.Cells.Item(13, 4).Value = 169160
$ws.Cells.Item(13, 5).Value = 30947
$ws.Cells.Item(13, 7).Value = 121
$ws.Cells.Item(13, 8).Value = 9863

# Row 27: Bielorrusia
$ws.Cells.Item(27, 2).Value = 59487
$ws.Cells.Item(27, 3).Value = 464
$ws.Cells.Item(27, 4).Value = 38688
$ws.Cells.Item(27, 5).Value = 20442
$ws.Cells.Item(27, 7).Value = 6
$ws.Cells.Item(27, 8).Value = 357

# Row 36: Kuwait
$ws.Cells.Item(36, 2).Value = 41033
$ws.Cells.Item(36, 3).Value = 742
$ws.Cells.Item(36, 4).Value = 32304
$ws.Cells.Item(36, 5).Value = 8395
$ws.Cells.Item(36, 7).Value = 4
$ws.Cells.Item(36, 8).Value = 334

# Row 42: was Suiza, now Filipinas (Filipinas overtakes Suiza in total cases)
$ws.Cells.Item(42, 1).Value = "Filipinas"
$ws.Cells.Item(42, 2).Value = 31825
$ws.Cells.Item(42, 3).Value = 1143
$ws.Cells.Item(42, 4).Value = 8442
$ws.Cells.Item(42, 5).Value = 22197
$ws.Cells.Item(42, 7).Value = 9
$ws.Cells.Item(42, 8).Value = 1186

# Row 43: was Filipinas, now Suiza
$ws.Cells.Item(43, 1).Value = "Suiza"
$ws.Cells.Item(43, 2).Value = 31332
$ws.Cells.Item(43, 3).Value = 22
$ws.Cells.Item(43, 4).Value = 29000
$ws.Cells.Item(43, 5).Value = 376
$ws.Cells.Item(43, 8).Value = 1956

# Row 69: Nepal
$ws.Cells.Item(69, 2).Value = 10099
$ws.Cells.Item(69, 3).Value = 538
$ws.Cells.Item(69, 4).Value = 2224
$ws.Cells.Item(69, 5).Value = 7851
$ws.Cells.Item(69, 7).Value = 1
$ws.Cells.Item(69, 8).Value = 24

# Row 77: Senegal
$ws.Cells.Item(77, 2).Value = 6034
$ws.Cells.Item(77, 3).Value = 64
$ws.Cells.Item(77, 4).Value = 4046
$ws.Cells.Item(77, 5).Value = 1899
$ws.Cells.Item(77, 7).Value = 3
$ws.Cells.Item(77, 8).Value = 89

# Row 78: Consejo Danes para los Refugiados
$ws.Cells.Item(78, 2).Value = 6027
$ws.Cells.Item(78, 3).Value = 103
$ws.Cells.Item(78, 4).Value = 861
$ws.Cells.Item(78, 5).Value = 5031

# Row 113: was Guinea Ecuatorial, now Madagascar (Madagascar overtakes Guinea Ecuatorial)
$ws.Cells.Item(113, 1).Value = "Madagascar"
$ws.Cells.Item(113, 2).Value = 1724
$ws.Cells.Item(113, 3).Value = 84
$ws.Cells.Item(113, 4).Value = 732
$ws.Cells.Item(113, 5).Value = 977
$ws.Cells.Item(113, 8).Value = 15

# Row 114: was Madagascar, now Guinea Ecuatorial
$ws.Cells.Item(114, 1).Value = "Guinea Ecuatorial"
$ws.Cells.Item(114, 2).Value = 1664
$ws.Cells.Item(114, 4).Value = 515
$ws.Cells.Item(114, 5).Value = 1117
$ws.Cells.Item(114, 8).Value = 32

# Row 134: Burkina Faso
$ws.Cells.Item(134, 2).Value = 907
$ws.Cells.Item(134, 3).Value = 4
$ws.Cells.Item(134, 4).Value = 823
$ws.Cells.Item(134, 5).Value = 31

# Row 138: Benin
$ws.Cells.Item(138, 2).Value = 850
$ws.Cells.Item(138, 3).Value = 43
$ws.Cells.Item(138, 4).Value = 272
$ws.Cells.Item(138, 5).Value = 565

# Row 145: San Marino
$ws.Cells.Item(145, 2).Value = 698
$ws.Cells.Item(145, 3).Value = 1
$ws.Cells.Item(145, 4).Value = 632
$ws.Cells.Item(145, 5).Value = 24
